$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 220.11111
$ws.Range("I9").Value = 113.5
$ws.Range("J9").Value = 433.33334
$ws.Range("K9").Value = 113.5
$ws.Range("L9").Value = 433.33334
$ws.Range("M9").Value = 55.5
$ws.Range("N9").Value = -771.33334

$ws.Range("H15").Value = 3217.0986
$ws.Range("I15").Value = 3217.0986
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 9651.2958
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -9482.2958

$ws.Range("H20").Value = 1999
$ws.Range("I20").Value = 1999
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1999
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1769

$ws.Range("H35").Value = 1999
$ws.Range("I35").Value = 1999
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1999
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1620

$ws.Range("H51").Value = 5527.273
$ws.Range("I51").Value = 2101
$ws.Range("J51").Value = 5869.9
$ws.Range("K51").Value = 2101
$ws.Range("L51").Value = 5869.9
$ws.Range("M51").Value = -1617
$ws.Range("N51").Value = -6837.9

$ws.Range("H82").Value = 231
$ws.Range("I82").Value = 231
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 693
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -287

$ws.Range("H85").Value = 231
$ws.Range("I85").Value = 231
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 693
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 4068.6667
$ws.Range("I14").Value = 4068.6667
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 4068.6667
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -3893.6667

$ws.Range("H86").Value = 34567
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 34567
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 34567
$ws.Range("N86").Value = -36939

$ws.Range("H89").Value = 34567
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 34567
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 103701
$ws.Range("N89").Value = -115557

$ws.Range("H132").Value = 66754.74000000001
$ws.Range("I132").Value = 49511.145
$ws.Range("J132").Value = 102966.3
$ws.Range("K132").Value = 148533.435
$ws.Range("L132").Value = 308898.9
$ws.Range("M132").Value = -146003.435
$ws.Range("N132").Value = -313958.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13885.637
$ws.Range("I86").Value = 14375.429
$ws.Range("J86").Value = 3600
$ws.Range("K86").Value = 14375.429
$ws.Range("L86").Value = 3600
$ws.Range("M86").Value = -13252.429
$ws.Range("N86").Value = -5846

$ws.Range("H89").Value = 13885.637
$ws.Range("I89").Value = 14375.429
$ws.Range("J89").Value = 3600
$ws.Range("K89").Value = 71877.145
$ws.Range("L89").Value = 18000
$ws.Range("M89").Value = -66261.145
$ws.Range("N89").Value = -29232

$ws.Range("H105").Value = 45456744
$ws.Range("I105").Value = 50002300
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 50002300
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = -50000553
$ws.Range("N105").Value = -4694

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12346943
$ws.Range("I16").Value = 1166.8334
$ws.Range("J16").Value = 37038496
$ws.Range("K16").Value = 1166.8334
$ws.Range("L16").Value = 37038496
$ws.Range("M16").Value = -879.8334
$ws.Range("N16").Value = -37039070

$ws.Range("H33").Value = 764.5
$ws.Range("I33").Value = 764.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 764.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -385.5

$ws.Range("H63").Value = 34271
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 34271
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 34271
$ws.Range("N63").Value = -35643

$ws.Range("H66").Value = 34271
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 34271
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 102813
$ws.Range("N66").Value = -109677

$ws.Range("H107").Value = 649.8570999999999
$ws.Range("I107").Value = 615.5
$ws.Range("J107").Value = 735.75
$ws.Range("K107").Value = 615.5
$ws.Range("L107").Value = 735.75
$ws.Range("M107").Value = 1304.5
$ws.Range("N107").Value = -4575.75

$ws.Range("H113").Value = 12346943
$ws.Range("I113").Value = 1166.8334
$ws.Range("J113").Value = 37038496
$ws.Range("K113").Value = 1166.8334
$ws.Range("L113").Value = 37038496
$ws.Range("M113").Value = 1003.1666
$ws.Range("N113").Value = -37042836

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 182.88889
$ws.Range("I38").Value = 73.666664
$ws.Range("J38").Value = 204.73334
$ws.Range("K38").Value = 220.999992
$ws.Range("L38").Value = 614.20002
$ws.Range("M38").Value = 126.000008
$ws.Range("N38").Value = -1308.20002

$ws.Range("H92").Value = 1250
$ws.Range("I92").Value = 1250
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 3750
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -2502
$ws.Range("N92").ClearContents()

$ws.Range("H131").Value = 880.2917
$ws.Range("I131").Value = 531.25
$ws.Range("J131").Value = 1054.8125
$ws.Range("K131").Value = 1593.75
$ws.Range("L131").Value = 3164.4375
$ws.Range("M131").Value = 3446.25
$ws.Range("N131").Value = -13244.4375

$ws.Range("H132").Value = 2816.7727
$ws.Range("I132").Value = 840.3333
$ws.Range("J132").Value = 3557.9375
$ws.Range("K132").Value = 7562.9997
$ws.Range("L132").Value = 32021.4375
$ws.Range("M132").Value = -5032.9997
$ws.Range("N132").Value = -37081.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3754.4707
$ws.Range("I80").Value = 3183.3333
$ws.Range("J80").Value = 4066
$ws.Range("K80").Value = 3183.3333
$ws.Range("L80").Value = 4066
$ws.Range("M80").Value = -2185.3333
$ws.Range("N80").Value = -6062

$ws.Range("H83").Value = 3754.4707
$ws.Range("I83").Value = 3183.3333
$ws.Range("J83").Value = 4066
$ws.Range("K83").Value = 15916.6665
$ws.Range("L83").Value = 20330
$ws.Range("M83").Value = -10924.6665
$ws.Range("N83").Value = -30314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1590.6
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1590.6
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1590.6
$ws.Range("N82").Value = -2312.6
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 1590.6
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1590.6
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1590.6
$ws.Range("N85").Value = -4086.6
$ws.Range("M85").ClearContents()

$ws.Range("H100").Value = 1589.65
$ws.Range("I100").Value = 1333.5834
$ws.Range("J100").Value = 1973.75
$ws.Range("K100").Value = 1333.5834
$ws.Range("L100").Value = 1973.75
$ws.Range("M100").Value = -792.5834

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
